# TD-6649 add Business Concept file manager domain name field
#
# The "term" sheet's header row labels column B "domain". Rename it to
# "domain_external_id" so the uploaded-hierarchy template exposes the
# domain's external id instead of its (ambiguous) display name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "domain_external_id"

# Leave the cursor on B2, matching where the edit was made.
$ws.Range("B2").Select()
